# Adds the 2021/12/13 (民國 110年12月13日) daily snapshot row to each of the
# five data sheets. Each sheet keeps its history in descending-date order
# starting at row 2, so the new day's numbers are inserted as a fresh row 2
# and all older rows shift down by one (handled natively by Rows.Insert()).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 台指期換倉成本計算 (A1:F19 -> A1:F20) ---------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2").Value = "日期：2021/12/13"
$ws1.Range("B2").Value = "202201"
$ws1.Range("C2").Value = 17691
$ws1.Range("D2").Value = 38344
$ws1.Range("E2").Value = 430864305
$ws1.Range("F2").Value = 17685

# --- Sheet 2: 散戶多空力道 (A1:B34 -> A1:B35) --------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = "日期：2021/12/13"
$ws2.Range("B2").Value = 0.06

# --- Sheet 3: 三大法人買賣金額 (A1:C34 -> A1:C35) ----------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
$ws3.Range("A2").Value = "110年12月13日"
$ws3.Range("B2").Value = -90.95
$ws3.Range("C2").Value = 32.29

# --- Sheet 4: 大盤多空點位 (A1:B33 -> A1:B34) --------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = "110年12月13日"
$ws4.Range("B2").Value = 17864.87

# --- Sheet 5: 期貨大額交易人未沖銷部位 (A1:N32 -> A1:N33) --------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
$ws5.Range("A2").Value = "2021/12/13"
$ws5.Range("B2").Value = 50961
$ws5.Range("C2").Value = 57854
$ws5.Range("D2").Value = 2566
$ws5.Range("E2").Value = 713
$ws5.Range("F2").Value = 27013
$ws5.Range("G2").Value = 52763
$ws5.Range("H2").Value = 2505
$ws5.Range("I2").Value = 1839
$ws5.Range("J2").Value = -25750
$ws5.Range("K2").Value = 666
$ws5.Range("L2").Value = 61
$ws5.Range("M2").Value = -1126
$ws5.Range("N2").Value = 1187
